$wb = $excel.ActiveWorkbook

# --- Sheet "Generic": NrBuckets (B4) goes from 3 to 4 (a 4th demand bucket,
#     index 3, is being introduced on the forecast sheets below) ---
$wsGeneric = $wb.Worksheets.Item("Generic")
$wsGeneric.Range("B4").Value = 4

# --- Sheet "Productdata": Leadtimes for Part_0001 (row4), Part_0002 (row5),
#     Part_0003 (row6) go from 445 to 890 ---
$wsProd = $wb.Worksheets.Item("Productdata")
$wsProd.Range("C4").Value = 890
$wsProd.Range("C5").Value = 890
$wsProd.Range("C6").Value = 890

# --- Sheet "ForecastedAverageDemand": add a new bucket row (row 5), matching
#     the formatting of the preceding bucket rows ---
$wsAvg = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvg.Range("A4").Copy()
$wsAvg.Range("A5").PasteSpecial(-4122)
$wsAvg.Range("A5").Value = 3
$wsAvg.Range("B5").Value = 0
$wsAvg.Range("C5").Value = 0
$wsAvg.Range("D5").Value = 0
$wsAvg.Range("E5").Value = 0
$wsAvg.Range("F5").Value = 0
$wsAvg.Range("G5").Value = 253
$wsAvg.Range("H5").Value = 45
$wsAvg.Range("I5").Value = 75

# --- Sheet "ForcastedStandardDeviation": add a new bucket row (row 5),
#     matching the formatting of the preceding bucket rows ---
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStd.Range("A4").Copy()
$wsStd.Range("A5").PasteSpecial(-4122)
$wsStd.Range("A5").Value = 3
$wsStd.Range("B5").Value = 0
$wsStd.Range("C5").Value = 0
$wsStd.Range("D5").Value = 0
$wsStd.Range("E5").Value = 0
$wsStd.Range("F5").Value = 0
$wsStd.Range("G5").Value = 36.62
$wsStd.Range("H5").Value = 1
$wsStd.Range("I5").Value = 2
